$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.02672833333333334
$ws.Range("H2").Value = 0.08018500000000001
$ws.Range("I2").Value = 0.3128877685602129
$ws.Range("J2").Value = 0.3128877685602129
$ws.Range("M2").Value = 13.71977066666667
$ws.Range("N2").Value = 41.159312
$ws.Range("O2").Value = 0.5515038136402627
$ws.Range("P2").Value = 0.5515038136402626
$ws.Range("Q2").Value = 0.3667066036355556
$ws.Range("R2").Value = 3.30035943272
$ws.Range("S2").Value = 0.1725587976023493
$ws.Range("T2").Value = 0.1725587976023493

$ws.Range("G3").Value = 0.02672833333333334
$ws.Range("H3").Value = 0.08018500000000001
$ws.Range("I3").Value = 0.3128877685602129
$ws.Range("J3").Value = 0.3128877685602129
$ws.Range("O3").Value = 0.172077867958883
$ws.Range("P3").Value = 0.1720778679588829
$ws.Range("Q3").Value = 0.1144182305894444
$ws.Range("R3").Value = 1.029764075305
$ws.Range("S3").Value = 0.05384106012425385
$ws.Range("T3").Value = 0.05384106012425385

$ws.Range("G4").Value = 0.02672833333333334
$ws.Range("H4").Value = 0.08018500000000001
$ws.Range("I4").Value = 0.3128877685602129
$ws.Range("J4").Value = 0.3128877685602129
$ws.Range("O4").Value = 0.2764183184008545
$ws.Range("P4").Value = 0.2764183184008545
$ws.Range("Q4").Value = 0.1837964130372222
$ws.Range("R4").Value = 1.654167717335
$ws.Range("S4").Value = 0.08648791083360979
$ws.Range("T4").Value = 0.08648791083360979

$ws.Range("G5").Value = 0.05869633333333333
$ws.Range("I5").Value = 0.6871122314397871
$ws.Range("J5").Value = 0.6871122314397871
$ws.Range("M5").Value = 13.71977066666667
$ws.Range("N5").Value = 41.159312
$ws.Range("O5").Value = 0.5515038136402627
$ws.Range("P5").Value = 0.5515038136402626
$ws.Range("Q5").Value = 0.8053002323075555
$ws.Range("R5").Value = 7.247702090768
$ws.Range("S5").Value = 0.3789450160379134
$ws.Range("T5").Value = 0.3789450160379133

$ws.Range("G6").Value = 0.05869633333333333
$ws.Range("I6").Value = 0.6871122314397871
$ws.Range("J6").Value = 0.6871122314397871
$ws.Range("O6").Value = 0.172077867958883
$ws.Range("P6").Value = 0.1720778679588829
$ws.Range("Q6").Value = 0.2512663441574444
$ws.Range("S6").Value = 0.1182368078346291
$ws.Range("T6").Value = 0.1182368078346291

$ws.Range("G7").Value = 0.05869633333333333
$ws.Range("I7").Value = 0.6871122314397871
$ws.Range("J7").Value = 0.6871122314397871
$ws.Range("O7").Value = 0.2764183184008545
$ws.Range("P7").Value = 0.2764183184008545
$ws.Range("S7").Value = 0.1899304075672447
$ws.Range("T7").Value = 0.1899304075672447

